$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision has been made to acquire rights for `"Barbie`" as the selected movie for Friday's assembly.`n"
$ws.Range("D2").Value = "Barbie_was_selected, "
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded.`n"
$ws.Range("D3").Value = "Barbie_was_selected, "
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" has been selected as the movie to show on Friday.`n"
$ws.Range("D4").Value = "Barbie_was_selected, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for the showing on Friday.`n"
$ws.Range("D5").Value = "Barbie_was_selected, "
$ws.Range("C6").Value = "MSG: None`n`nMSG: I have acquired the rights to both movies.`n"
$ws.Range("D6").Value = "both_movies, "
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision process resulted in no agreement on which movie to show on Friday.`n"
$ws.Range("D7").Value = "no_decision, "
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding which movie will be shown on Friday.`n"
$ws.Range("D8").Value = "no_decision, "
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie will be selected for Friday.`n"
$ws.Range("D10").Value = "no_decision, "
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has resulted in no agreement, so no movie rights will be acquired at this time.`n"
$ws.Range("D11").Value = "no_decision, "
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D12").Value = "no_decision, "
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("D13").Value = "Oppenheimer_was_selected, "
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision was made that there would be no movie selected for Friday.`n"
$ws.Range("D14").Value = "no_decision, "
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("D15").Value = "Barbie_was_selected, "
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision about the movie to be shown on Friday resulted in no agreement, so no decision was made.`n"
$ws.Range("D16").Value = "no_decision, "
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no agreement was reached about which movie to show on Friday.`n"
$ws.Range("D17").Value = "no_decision, "
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been recorded successfully.`n"
$ws.Range("D18").Value = "both_movies, "
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for showing on Friday.`n"
$ws.Range("D19").Value = "Barbie_was_selected, "
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" has been selected for the assembly.`n"
$ws.Range("D20").Value = "Barbie_was_selected, "
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision process did not result in a clear choice, and thus no movie has been selected for Friday.`n"
$ws.Range("D21").Value = "no_decision, "
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday. If further discussions are needed, they can take place at a later time.`n"
$ws.Range("D22").Value = "no_decision, "
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for both movies.`n"
$ws.Range("D23").Value = "both_movies, "
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been made.`n"
$ws.Range("D24").Value = "Barbie_was_selected, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for both movies.`n"
$ws.Range("D25").Value = "both_movies, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision-making process did not lead to a consensus on which movie to show on Friday, resulting in no decision being made.`n"
$ws.Range("D26").Value = "no_decision, "
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision being made about which movie to show on Friday.`n"
$ws.Range("D27").Value = "no_decision, "
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been made.`n"
$ws.Range("D28").Value = "both_movies, "
$ws.Range("C29").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("D29").Value = "both_movies, "
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been successfully recorded.`n"
$ws.Range("D30").Value = "both_movies, "
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has ended without a conclusive choice.`n"
$ws.Range("D31").Value = "no_decision, "
$ws.Range("C32").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("D32").Value = "both_movies, "
$ws.Range("C33").Value = "MSG: None`n`nMSG: The rights to `"Barbie`" have been acquired for showing on Friday.`n"
$ws.Range("D33").Value = "Barbie_was_selected, "
$ws.Range("C34").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie.`"`n"
$ws.Range("D34").Value = "Barbie_was_selected, "
$ws.Range("C35").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("D35").Value = "both_movies, "
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D36").Value = "no_decision, "
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded successfully.`n"
$ws.Range("D37").Value = "both_movies, "
$ws.Range("C38").Value = "MSG: None`n`nMSG: The rights for both movies, `"Oppenheimer`" and `"Barbie,`" have been acquired.`n"
$ws.Range("D38").Value = "both_movies, "
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Oppenheimer.`"`n"
$ws.Range("D39").Value = "Oppenheimer_was_selected, "
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded: no agreement was reached regarding the movie to be shown on Friday.`n"
$ws.Range("D40").Value = "no_decision, "
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision to acquire the rights for the movie `"Barbie`" has been recorded successfully.`n"
$ws.Range("D41").Value = "Barbie_was_selected, "
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been finalized, so I have recorded the outcome as no decision being made.`n"
$ws.Range("D42").Value = "no_decision, "
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision process has concluded without selecting a movie for Friday.`n"
$ws.Range("D43").Value = "no_decision, "
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision has been recorded, and no specific movie was chosen for Friday.`n"
$ws.Range("D44").Value = "no_decision, "
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded, and the conclusion reached is that no movie has been selected for Friday.`n"
$ws.Range("D45").Value = "no_decision, "
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision process concluded without a definitive choice for Friday's movie. Therefore, I have recorded the decision using the no_decision function.`n"
$ws.Range("D46").Value = "no_decision, "
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been recorded successfully.`n"
$ws.Range("D47").Value = "Oppenheimer_was_selected, "
$ws.Range("C48").Value = "MSG: None`n`nMSG: The selection for the movie `"Barbie`" has been successfully recorded.`n"
$ws.Range("D48").Value = "Barbie_was_selected, "
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision about which movie to play on Friday did not reach a consensus among the committee members, resulting in no decision being made.`n"
$ws.Range("D49").Value = "no_decision, "
$ws.Range("C50").Value = "MSG: None`n`nMSG: The committee ended the conversation without reaching a decision about which movie will be shown on Friday.`n"
$ws.Range("D50").Value = "no_decision, "
$ws.Range("C51").Value = "MSG: None`n`nMSG: None`n`nMSG: The committee did not reach a decision on which movie to show on Friday, so I will proceed to call the no_decision function.`n"
$ws.Range("D51").Value = "no_decision, , no_decision, "
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision was to not select any movie for Friday.`n"
$ws.Range("D52").Value = "no_decision, "
$ws.Range("C53").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire the rights for both `"Oppenheimer`" and `"Barbie`" for Friday's screening.`n"
$ws.Range("D53").Value = "both_movies, "
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both movies: `"Barbie`" and `"Oppenheimer.`"`n"
$ws.Range("D54").Value = "both_movies, "
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision process ended without a clear selection for Friday's movie, so no movie rights will be acquired at this time.`n"
$ws.Range("D55").Value = "no_decision, "
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision has been recorded, and the outcome is that no movie was selected for the Friday assembly.`n"
$ws.Range("D56").Value = "no_decision, "
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("D57").Value = "Barbie_was_selected, "
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision to acquire a movie for Friday was not reached by the committee, resulting in a no-decision outcome.`n"
$ws.Range("D58").Value = "no_decision, "
$ws.Range("C59").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("D59").Value = "Barbie_was_selected, "
$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has been recorded as no decision was reached.`n"
$ws.Range("D60").Value = "no_decision, "
$ws.Range("C61").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for both movies.`n"
$ws.Range("D61").Value = "both_movies, "
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday was made.`n"
$ws.Range("D62").Value = "no_decision, "
$ws.Range("C63").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("D63").Value = "both_movies, "
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("D64").Value = "Barbie_was_selected, "
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded successfully.`n"
$ws.Range("D65").Value = "both_movies, "
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Oppenheimer`" will be acquired for the showing on Friday.`n"
$ws.Range("D66").Value = "Oppenheimer_was_selected, "
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision process concluded without a selection for Friday's movie.`n"
$ws.Range("D67").Value = "no_decision, "
$ws.Range("C68").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for `"Barbie`" will be acquired for the upcoming assembly.`n"
$ws.Range("D68").Value = "Barbie_was_selected, "
$ws.Range("C69").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("D69").Value = "no_decision, "
